$d = $word.ActiveDocument

# 1. Title: "CampuSales" -> "CampuSale"
$d.Content.Find.Execute("CampuSales", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CampuSale", 2) | Out-Null

# 2. "Added referential integrity to database." -> "Reconstructed database to add referential integrity."
$d.Content.Find.Execute("Added referential integrity to database.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Reconstructed database to add referential integrity.", 2) | Out-Null

# 3. Insert two new bullet items before "Link items on the listing page to their specific item pages."
#    (which keeps them ahead of that item, right after "Now that the front end ... post an item for sale.")
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Link items on the listing page")) {
        $target = $i
        break
    }
}

# Insert a new paragraph before the target, splitting off its formatting (numbered-list style),
# then fill in the text for "Display test items from database on listing page."
$d.Paragraphs.Item($target).Range.InsertParagraphBefore()
$d.Paragraphs.Item($target).Range.Text = "Display test items from database on listing page."

# Insert a second new paragraph before the (still) target paragraph for
# "Make dynamic route for item pages."
$d.Paragraphs.Item($target + 1).Range.InsertParagraphBefore()
$d.Paragraphs.Item($target + 1).Range.Text = "Make dynamic route for item pages."

# 4. Remove the bullet "Link info on profile page to the database."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Link info on profile page to the database")) {
        $d.Paragraphs.Item($i).Range.Delete()
        break
    }
}

Write-Output "done"
